$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clauses = @("1.1","1.2","1.3","1.4","1.5","1.6","1.7","1.8","2.1","2.2","2.3","2.4","2.5","2.6","2.7","2.8","3.1","3.2","3.3","3.4","3.5","4.1","4.2","4.3","4.4")

for ($i = 0; $i -lt $clauses.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = "MVSP-" + $clauses[$i]
}

$ws.Range("E2").Select()
